$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value would otherwise be auto-parsed as a number;
# force Text format so the stored type matches the source (inline string).
$textCells = @("D5","D6","D14","D19","D20","D22","D25","D27","D32","D38","D39","D44","D46","D48","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = "68.333.67"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "2.640.55"
$ws.Range("E3").Value = "  +0.03%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "597.20"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").Value = "154.45"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E8").Value = "  -0.66%  "
$ws.Range("D9").Value = "2.639.65"
$ws.Range("E10").Value = "  +6.83%  "
$ws.Range("E11").Value = "  -0.61%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "28.08"
$ws.Range("E14").Value = "  +1.50%  "
$ws.Range("E15").Value = "  +1.93%  "
$ws.Range("D16").Value = "3.119.81"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "68.238.94"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "2.620.57"
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").Value = "11.39"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "363.71"
$ws.Range("E20").Value = "  -2.68%  "
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").Value = "4.37"
$ws.Range("E22").Value = "  +2.70%  "
$ws.Range("E23").Value = "  +1.31%  "
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "74.64"
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "9.83"
$ws.Range("E27").Value = "  -1.42%  "
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("D29").Value = "2.778.64"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.23%  "
$ws.Range("E31").Value = "  -0.71%  "
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  +3.34%  "
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  +1.23%  "
$ws.Range("E35").Value = "  +3.45%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("D38").Value = "160.38"
$ws.Range("E38").Value = "  +1.42%  "
$ws.Range("D39").Value = "19.36"
$ws.Range("E39").Value = "  +1.00%  "
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  +0.88%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₆0336"
$ws.Range("E43").Value = "  +4.91%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.66"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D46").Value = "40.68"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").Value = "156.28"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +1.80%  "
$ws.Range("D50").Value = "1.71"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "21.81"
$ws.Range("E51").Value = "  -1.25%  "

# Restore default (Normal) style on the forced-text cells so only the value changed
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
